# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions" (Thu Nov  2 06:14:19 UTC 2023).
# All target cells are plain text (prices/percentages are formatted strings,
# not numeric cell values), so cells whose new text would otherwise be
# auto-coerced into a Number by the COM Value setter are first forced to the
# Text number format ("@") to preserve e.g. "1.00" / "0.0520" verbatim.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.360.32"
$ws.Range("E2").Value = "  +2.69%  "
$ws.Range("D3").Value = "1.846.85"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.64"
$ws.Range("E5").Value = "  +2.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("E6").Value = "  +4.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.32"
$ws.Range("E8").Value = "  +11.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.305"
$ws.Range("E9").Value = "  +6.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0692"
$ws.Range("E10").Value = "  +3.37%  "
$ws.Range("E11").Value = "  +3.59%  "
$ws.Range("D12").Value = "2.113.67"
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.34"
$ws.Range("E13").Value = "  +2.77%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.835.26"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.668"
$ws.Range("E15").Value = "  +6.85%  "
$ws.Range("E16").Value = "  +6.72%  "
$ws.Range("D17").Value = "35.368.58"
$ws.Range("E17").Value = "  +2.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.69"
$ws.Range("E18").Value = "  +4.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.97"
$ws.Range("E19").Value = "  +2.32%  "
$ws.Range("D20").Value = "0.0₃0798"
$ws.Range("E20").Value = "  +4.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.03"
$ws.Range("E21").Value = "  +9.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.61"
$ws.Range("E22").Value = "  +13.15%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.04"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("E26").Value = "  +2.79%  "
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.38"
$ws.Range("E29").Value = "  +12.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").Value = "3.322.49"
$ws.Range("E31").Value = "  +36.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0543"
$ws.Range("E32").Value = "  +6.44%  "
$ws.Range("E33").Value = "  +4.95%  "
$ws.Range("E34").Value = "  +5.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.85"
$ws.Range("E35").Value = "  +2.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "95.92"
$ws.Range("E36").Value = "  +17.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.682"
$ws.Range("E37").Value = "  +7.30%  "
$ws.Range("D38").Value = "1.350.95"
$ws.Range("E38").Value = "  +2.17%  "
$ws.Range("E39").Value = "  +3.35%  "
$ws.Range("E40").Value = "  +5.86%  "
$ws.Range("E41").Value = "  +3.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  +6.28%  "
$ws.Range("E43").Value = "  +4.16%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.63"
$ws.Range("E44").Value = "  +7.45%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0520"
$ws.Range("E47").Value = "  +2.00%  "
$ws.Range("E48").Value = "  +8.27%  "
$ws.Range("D49").Value = "2.013.18"
$ws.Range("E49").Value = "  +2.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.39"
$ws.Range("E51").Value = "  +1.59%  "
